$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply weekly re-dating / re-ordering of Frutilla price rows (Agricola del Norte S.A. de Arica).
# Each target row below receives the full data set (Fecha, Calidad, Volumen, Precios) that
# previously belonged to a different row, per the authoritative diff.
# Row 2 <= data previously in row 32
$ws.Range("D2").Value2 = 44249
$ws.Range("L2").Value2 = "Especial"
$ws.Range("M2").Value2 = 200
$ws.Range("N2").Value2 = 6000
$ws.Range("O2").Value2 = 7000
$ws.Range("P2").Value2 = 6500
$ws.Range("S2").Value2 = 2167

# Row 3 <= data previously in row 33
$ws.Range("D3").Value2 = 44249
$ws.Range("L3").Value2 = "Primera"
$ws.Range("M3").Value2 = 160
$ws.Range("N3").Value2 = 4500
$ws.Range("O3").Value2 = 5000
$ws.Range("P3").Value2 = 4750
$ws.Range("S3").Value2 = 1583

# Row 4 <= data previously in row 6
$ws.Range("D4").Value2 = 44242
$ws.Range("L4").Value2 = "Especial"
$ws.Range("M4").Value2 = 50
$ws.Range("N4").Value2 = 7000
$ws.Range("O4").Value2 = 8000
$ws.Range("P4").Value2 = 7500
$ws.Range("S4").Value2 = 2500

# Row 5 <= data previously in row 7
$ws.Range("D5").Value2 = 44242
$ws.Range("L5").Value2 = "Primera"
$ws.Range("M5").Value2 = 90
$ws.Range("N5").Value2 = 6000
$ws.Range("O5").Value2 = 7000
$ws.Range("P5").Value2 = 6500
$ws.Range("S5").Value2 = 2167

# Row 6 <= data previously in row 8
$ws.Range("D6").Value2 = 44242
$ws.Range("L6").Value2 = "Segunda"
$ws.Range("M6").Value2 = 100
$ws.Range("N6").Value2 = 4000
$ws.Range("O6").Value2 = 5000
$ws.Range("P6").Value2 = 4500
$ws.Range("S6").Value2 = 1500

# Row 7 <= data previously in row 18
$ws.Range("D7").Value2 = 44322
$ws.Range("L7").Value2 = "Especial"
$ws.Range("M7").Value2 = 200
$ws.Range("N7").Value2 = 7000
$ws.Range("O7").Value2 = 7500
$ws.Range("P7").Value2 = 7250
$ws.Range("S7").Value2 = 2417

# Row 8 <= data previously in row 19
$ws.Range("D8").Value2 = 44322
$ws.Range("L8").Value2 = "Primera"
$ws.Range("M8").Value2 = 160
$ws.Range("N8").Value2 = 6000
$ws.Range("O8").Value2 = 6500
$ws.Range("P8").Value2 = 6250
$ws.Range("S8").Value2 = 2083

# Row 9 <= data previously in row 20
$ws.Range("D9").Value2 = 44322
$ws.Range("L9").Value2 = "Segunda"
$ws.Range("M9").Value2 = 100
$ws.Range("N9").Value2 = 5000
$ws.Range("O9").Value2 = 5500
$ws.Range("P9").Value2 = 5250
$ws.Range("S9").Value2 = 1750

# Row 10 <= data previously in row 12
$ws.Range("D10").Value2 = 44351
$ws.Range("L10").Value2 = "Especial"
$ws.Range("M10").Value2 = 160
$ws.Range("N10").Value2 = 7500
$ws.Range("O10").Value2 = 8000
$ws.Range("P10").Value2 = 7750
$ws.Range("S10").Value2 = 2583

# Row 11 <= data previously in row 13
$ws.Range("D11").Value2 = 44351
$ws.Range("L11").Value2 = "Primera"
$ws.Range("M11").Value2 = 100
$ws.Range("N11").Value2 = 6000
$ws.Range("O11").Value2 = 6500
$ws.Range("P11").Value2 = 6250
$ws.Range("S11").Value2 = 2083

# Row 12 <= data previously in row 14
$ws.Range("D12").Value2 = 44351
$ws.Range("L12").Value2 = "Segunda"
$ws.Range("M12").Value2 = 200
$ws.Range("N12").Value2 = 4500
$ws.Range("O12").Value2 = 5000
$ws.Range("P12").Value2 = 4750
$ws.Range("S12").Value2 = 1583

# Row 13 <= data previously in row 2
$ws.Range("D13").Value2 = 44172
$ws.Range("L13").Value2 = "Especial"
$ws.Range("M13").Value2 = 100
$ws.Range("N13").Value2 = 6500
$ws.Range("O13").Value2 = 7000
$ws.Range("P13").Value2 = 6750
$ws.Range("S13").Value2 = 2250

# Row 14 <= data previously in row 3
$ws.Range("D14").Value2 = 44172
$ws.Range("L14").Value2 = "Primera"
$ws.Range("M14").Value2 = 160
$ws.Range("N14").Value2 = 5500
$ws.Range("O14").Value2 = 6000
$ws.Range("P14").Value2 = 5750
$ws.Range("S14").Value2 = 1917

# Row 15 <= data previously in row 4
$ws.Range("D15").Value2 = 44172
$ws.Range("L15").Value2 = "Segunda"
$ws.Range("M15").Value2 = 160
$ws.Range("N15").Value2 = 5000
$ws.Range("O15").Value2 = 5500
$ws.Range("P15").Value2 = 5250
$ws.Range("S15").Value2 = 1750

# Row 16 <= data previously in row 5
$ws.Range("D16").Value2 = 44172
$ws.Range("L16").Value2 = "Tercera"
$ws.Range("M16").Value2 = 140
$ws.Range("N16").Value2 = 3500
$ws.Range("O16").Value2 = 4000
$ws.Range("P16").Value2 = 3750
$ws.Range("S16").Value2 = 1250

# Row 17 <= data previously in row 25
$ws.Range("D17").Value2 = 44596
$ws.Range("L17").Value2 = "Especial"
$ws.Range("M17").Value2 = 100
$ws.Range("N17").Value2 = 8000
$ws.Range("O17").Value2 = 9000
$ws.Range("P17").Value2 = 8500
$ws.Range("S17").Value2 = 2833

# Row 18 <= data previously in row 26
$ws.Range("D18").Value2 = 44596
$ws.Range("L18").Value2 = "Primera"
$ws.Range("M18").Value2 = 130
$ws.Range("N18").Value2 = 6000
$ws.Range("O18").Value2 = 7000
$ws.Range("P18").Value2 = 6500
$ws.Range("S18").Value2 = 2167

# Row 19 <= data previously in row 27
$ws.Range("D19").Value2 = 44596
$ws.Range("L19").Value2 = "Segunda"
$ws.Range("M19").Value2 = 160
$ws.Range("N19").Value2 = 5000
$ws.Range("O19").Value2 = 6000
$ws.Range("P19").Value2 = 5500
$ws.Range("S19").Value2 = 1833

# Row 20 <= data previously in row 28
$ws.Range("D20").Value2 = 44596
$ws.Range("L20").Value2 = "Tercera"
$ws.Range("M20").Value2 = 100
$ws.Range("N20").Value2 = 4000
$ws.Range("O20").Value2 = 5000
$ws.Range("P20").Value2 = 4500
$ws.Range("S20").Value2 = 1500

# Row 25 <= data previously in row 29
$ws.Range("D25").Value2 = 44708
$ws.Range("L25").Value2 = "Primera"
$ws.Range("M25").Value2 = 50
$ws.Range("N25").Value2 = 6000
$ws.Range("O25").Value2 = 7000
$ws.Range("P25").Value2 = 6500
$ws.Range("S25").Value2 = 2167

# Row 26 <= data previously in row 30
$ws.Range("D26").Value2 = 44708
$ws.Range("L26").Value2 = "Segunda"
$ws.Range("M26").Value2 = 60
$ws.Range("N26").Value2 = 4000
$ws.Range("O26").Value2 = 5000
$ws.Range("P26").Value2 = 4500
$ws.Range("S26").Value2 = 1500

# Row 27 <= data previously in row 31
$ws.Range("D27").Value2 = 44708
$ws.Range("L27").Value2 = "Tercera"
$ws.Range("M27").Value2 = 50
$ws.Range("N27").Value2 = 3000
$ws.Range("O27").Value2 = 4000
$ws.Range("P27").Value2 = 3500
$ws.Range("S27").Value2 = 1167

# Row 28 <= data previously in row 9
$ws.Range("D28").Value2 = 44389
$ws.Range("L28").Value2 = "Especial"
$ws.Range("M28").Value2 = 100
$ws.Range("N28").Value2 = 7500
$ws.Range("O28").Value2 = 8000
$ws.Range("P28").Value2 = 7750
$ws.Range("S28").Value2 = 2583

# Row 29 <= data previously in row 10
$ws.Range("D29").Value2 = 44389
$ws.Range("L29").Value2 = "Primera"
$ws.Range("M29").Value2 = 160
$ws.Range("N29").Value2 = 6000
$ws.Range("O29").Value2 = 7000
$ws.Range("P29").Value2 = 6500
$ws.Range("S29").Value2 = 2167

# Row 30 <= data previously in row 11
$ws.Range("D30").Value2 = 44389
$ws.Range("L30").Value2 = "Segunda"
$ws.Range("M30").Value2 = 200
$ws.Range("N30").Value2 = 5500
$ws.Range("O30").Value2 = 6000
$ws.Range("P30").Value2 = 5750
$ws.Range("S30").Value2 = 1917

# Row 31 <= data previously in row 15
$ws.Range("D31").Value2 = 44200
$ws.Range("L31").Value2 = "Especial"
$ws.Range("M31").Value2 = 50
$ws.Range("N31").Value2 = 4500
$ws.Range("O31").Value2 = 5000
$ws.Range("P31").Value2 = 4750
$ws.Range("S31").Value2 = 1583

# Row 32 <= data previously in row 16
$ws.Range("D32").Value2 = 44200
$ws.Range("L32").Value2 = "Primera"
$ws.Range("M32").Value2 = 80
$ws.Range("N32").Value2 = 3500
$ws.Range("O32").Value2 = 4000
$ws.Range("P32").Value2 = 3750
$ws.Range("S32").Value2 = 1250

# Row 33 <= data previously in row 17
$ws.Range("D33").Value2 = 44200
$ws.Range("L33").Value2 = "Segunda"
$ws.Range("M33").Value2 = 120
$ws.Range("N33").Value2 = 2500
$ws.Range("O33").Value2 = 3000
$ws.Range("P33").Value2 = 2750
$ws.Range("S33").Value2 = 917
